$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Step 1: update row 6 (currently "2.1.5" row) to become the new "2.1.7" row,
#     including an extra bullet paragraph with the bug-fix note. ---
$row6 = $t.Rows.Item(6)
$row6.Cells.Item(1).Range.Text = "2.1.7"
$row6.Cells.Item(2).Range.Text = "6/12/2019"
$row6.Cells.Item(3).Range.Text = "Able to do 2D STRS`rFix the bug of binomial process for enrichment"

# --- Step 2: update row 5 (currently "2.1.3" row) to become the new "2.1.5" row. ---
$row5 = $t.Rows.Item(5)
$row5.Cells.Item(1).Range.Text = "2.1.5"
$row5.Cells.Item(2).Range.Text = "5/13/2019"
$row5.Cells.Item(3).Range.Text = "Able to visualize binary results for enrichment"

# --- Step 3: delete the old "2.1.7" row (now row 7), whose content has been folded
#     into row 6 above. ---
$t2 = $d.Tables.Item(1)
$row7 = $t2.Rows.Item(7)
$row7.Delete()

# --- Step 4: move the "_GoBack" bookmark from the end of the document (after the
#     last picture) to the start of the (new) "2.1.5" row's version-number cell. ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldBm = $d.Bookmarks.Item("_GoBack")
    $oldBm.Delete()
}

$t3 = $d.Tables.Item(1)
$verCell = $t3.Rows.Item(5).Cells.Item(1)
$startPos = $verCell.Range.Start
$ip = $d.Range($startPos, $startPos)
$d.Bookmarks.Add("_GoBack", $ip)

Write-Output "done"
